$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 279, pushing the
# existing rows 279-292 down to 281-294 (matches dimension change to A1:T294).
$ws.Range("A279:A280").EntireRow.Insert()

# New row 279: Naranja / Lane Late / Primera, fecha 2021-11-09 (serial 44509)
$ws.Range("A279").Value = 7
$ws.Range("B279").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C279").Value = "Ñuble"
$ws.Range("D279").Value2 = 44509
$ws.Range("E279").Value = 16
$ws.Range("F279").Value = "Fruta"
$ws.Range("G279").Value = 100102
$ws.Range("H279").Value = "Cítricos"
$ws.Range("I279").Value = 100102005
$ws.Range("J279").Value = "Naranja"
$ws.Range("K279").Value = "Lane Late"
$ws.Range("L279").Value = "Primera"
$ws.Range("M279").Value = 160
$ws.Range("N279").Value = 8000
$ws.Range("O279").Value = 8500
$ws.Range("P279").Value = 8250
$ws.Range("Q279").Value = "$/bandeja 15 kilos granel"
$ws.Range("R279").Value = "Región de O'Higgins"
$ws.Range("S279").Value = 550
$ws.Range("T279").Value = 15

# New row 280: Naranja / Lane Late / Segunda, fecha 2021-11-09 (serial 44509)
$ws.Range("A280").Value = 7
$ws.Range("B280").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C280").Value = "Ñuble"
$ws.Range("D280").Value2 = 44509
$ws.Range("E280").Value = 16
$ws.Range("F280").Value = "Fruta"
$ws.Range("G280").Value = 100102
$ws.Range("H280").Value = "Cítricos"
$ws.Range("I280").Value = 100102005
$ws.Range("J280").Value = "Naranja"
$ws.Range("K280").Value = "Lane Late"
$ws.Range("L280").Value = "Segunda"
$ws.Range("M280").Value = 120
$ws.Range("N280").Value = 7000
$ws.Range("O280").Value = 7500
$ws.Range("P280").Value = 7250
$ws.Range("Q280").Value = "$/bandeja 15 kilos granel"
$ws.Range("R280").Value = "Región de O'Higgins"
$ws.Range("S280").Value = 483
$ws.Range("T280").Value = 15
